$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the ID column values down (A4 -> A5 -> A6 -> A7), and set the new A4 value
$ws.Range("A6").Value = "q1371623"
$ws.Range("A5").Value = "q0762379"
$ws.Range("A4").Value = "c1243957"

# Add the new row 7 with data (previous A6 id moves here along with fresh S1/S2/S3/Q1/Q2/Q3 values)
$ws.Range("A7").Value = "q1411379"
$ws.Range("B7").Value = 1.43264
$ws.Range("C7").Value = 0.32398
$ws.Range("D7").Value = 0.84935
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = 8
